$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously last data row (row 45) represented a record dated 2021-10-22
# (serial 44491) with volume 5000 and prices 850. A weekly update inserts a
# new record (dated 2021-11-09 / serial 44509, volume 6000, prices 800) as
# the new row 45, and pushes the old record down to the new row 46.

# 1) Create new row 46 with the values that used to live in row 45.
$ws.Range("A46").Value = $ws.Range("A45").Value()
$ws.Range("B46").Value = $ws.Range("B45").Value()
$ws.Range("C46").Value = $ws.Range("C45").Value()

$ws.Range("D46").Value = 44491
$ws.Range("D46").NumberFormat = $ws.Range("D45").NumberFormat()

$ws.Range("E46").Value = $ws.Range("E45").Value()
$ws.Range("F46").Value = $ws.Range("F45").Value()
$ws.Range("G46").Value = $ws.Range("G45").Value()
$ws.Range("H46").Value = $ws.Range("H45").Value()
$ws.Range("I46").Value = $ws.Range("I45").Value()

$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 850
$ws.Range("L46").Value = 850
$ws.Range("M46").Value = 850

$ws.Range("N46").Value = $ws.Range("N45").Value()
$ws.Range("O46").Value = $ws.Range("O45").Value()

$ws.Range("P46").Value = 850
$ws.Range("Q46").Value = 1

$ws.Range("R46").Value = $ws.Range("R45").Value()

# 2) Overwrite row 45 with the new record's values (columns A,B,C,E,F,G,H,I,N,O,Q,R stay the same).
$ws.Range("D45").Value = 44509
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 800
$ws.Range("L45").Value = 800
$ws.Range("M45").Value = 800
$ws.Range("P45").Value = 800
